# ---------------------------------------------------------------------------
# Applies the "excel data source files udpates" commit:
#   1. pitstop sheet: two new trailing columns (tyre_before / tyre_after)
#      with per-row tyre-compound data.
#   2. Two brand-new sheets appended at the end of the workbook:
#        - weather  : key/value weather readings for the race
#        - altitude : a single delta reading
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. pitstop sheet -> add tyre_before / tyre_after columns
# ---------------------------------------------------------------------------
$pitstop = $wb.Worksheets.Item(4)

$pitstop.Range("H1").Value = "tyre_before"
$pitstop.Range("I1").Value = "tyre_after"

# Row (2..38) -> [tyre_before, tyre_after]
$tyreData = @(
    @(4,2), @(3,2), @(2,3), @(3,4), @(3,2), @(2,3), @(3,2), @(4,2), @(2,3),
    @(4,2), @(2,3), @(3,4), @(3,2), @(3,2), @(2,3), @(3,2), @(2,4), @(4,2),
    @(2,3), @(4,2), @(3,3), @(3,2), @(2,4), @(4,2), @(3,2), @(2,3), @(3,2),
    @(2,3), @(3,2), @(2,4), @(3,2), @(2,4), @(3,2), @(2,3), @(3,2), @(4,3),
    @(3,3)
)

for ($i = 0; $i -lt $tyreData.Count; $i++) {
    $r = $i + 2
    $pitstop.Cells.Item($r, 8).Value = $tyreData[$i][0]
    $pitstop.Cells.Item($r, 9).Value = $tyreData[$i][1]
}

$pitstop.Columns.Item(8).AutoFit() | Out-Null
$pitstop.Columns.Item(9).AutoFit() | Out-Null

$pitstop.Range("H39").Select()

# ---------------------------------------------------------------------------
# 2. Append the two new sheets at the end: weather, altitude
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$weather = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$weather.Name = "weather"

$altitude = $wb.Worksheets.Add([System.Type]::Missing, $weather)
$altitude.Name = "altitude"

# --- weather sheet content ---
$weather.Range("A1").Value = "Skycondition"
$weather.Range("B1").Value = "Breezy and Overcast"

$weather.Range("A2").Value = "Temperature"
$weather.Range("B2").Value = "66.04°F"

$weather.Range("A3").Value = "Humidity"
$weather.Range("B3").Value = 0.45
$weather.Range("B3").NumberFormat = "0%"

$weather.Range("A4").Value = "Wind speed"
$weather.Range("B4").Value = "15.44 mph"

$weather.Range("A5").Value = "Wind bearing"
$weather.Range("B5").Value = "33°"

$weather.Columns.Item(1).AutoFit() | Out-Null
$weather.Columns.Item(2).AutoFit() | Out-Null

$weather.PageSetup.Orientation = 1

$weather.Range("D3").Select()

# --- altitude sheet content ---
$altitude.Range("A1").Value = "delta"
$altitude.Range("B1").Value = 7.4

# Last selection made wins for "ActiveSheet" / tabSelected bookkeeping, so
# do this last -> altitude ends up the active sheet, matching the target.
$altitude.Range("B2").Select()
